$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FindCarTest")

$ws.Range("D1").Value = "carTitle"
$ws.Range("D2").Value = "BMW Cars"
$ws.Range("D4").Value = "Toyota Cars"
$ws.Range("D3").Value = "Maruti Cars"

$ws.Range("D4").Select()
